$d = $word.ActiveDocument

$replacements = @(
    @("811×9=", "675×7="),
    @("576×8=", "219×4="),
    @("411×9=", "778×9="),
    @("159×2=", "778×4="),
    @("515×7=", "749×6="),
    @("379×7=", "848×4="),
    @("411×8=", "429×7="),
    @("275×8=", "647×2="),
    @("778×6=", "189×2="),
    @("157×5=", "155×6="),
    @("573×9=", "737×6="),
    @("806×6=", "970×3="),
    @("987×5=", "981×9="),
    @("591×7=", "942×6="),
    @("142×4=", "965×5="),
    @("434×4=", "843×7="),
    @("850×7=", "402×9="),
    @("149×7=", "417×7="),
    @("165×8=", "242×2="),
    @("212×9=", "169×3="),
    @("331×7=", "547×6="),
    @("864×5=", "949×3="),
    @("616×9=", "137×2="),
    @("345×5=", "628×3="),
    @("943×9=", "309×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
